$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: shift existing Mandarina price rows down by two and add the two newest
# (289-290) at the top and two more (339-340) at the bottom, matching upstream report pull.

$ws.Range("D289").Value = 45124
$ws.Range("K289").Value = 'Clementina'
$ws.Range("L289").Value = 'Primera'
$ws.Range("M289").Value = 80
$ws.Range("N289").Value = 10000
$ws.Range("O289").Value = 10000
$ws.Range("P289").Value = 10000
$ws.Range("Q289").Value = '$/bandeja 10 kilos'
$ws.Range("R289").Value = 'Región de O''Higgins'
$ws.Range("S289").Value = 1000
$ws.Range("T289").Value = 10

$ws.Range("D290").Value = 45124
$ws.Range("K290").Value = 'Clementina'
$ws.Range("L290").Value = 'Segunda'
$ws.Range("M290").Value = 80
$ws.Range("N290").Value = 8000
$ws.Range("O290").Value = 8000
$ws.Range("P290").Value = 8000
$ws.Range("Q290").Value = '$/bandeja 10 kilos'
$ws.Range("R290").Value = 'Región de O''Higgins'
$ws.Range("S290").Value = 800
$ws.Range("T290").Value = 10

$ws.Range("D291").Value = 45070
$ws.Range("K291").Value = 'Clementina'
$ws.Range("L291").Value = 'Especial'
$ws.Range("M291").Value = 80
$ws.Range("N291").Value = 14000
$ws.Range("O291").Value = 14000
$ws.Range("P291").Value = 14000
$ws.Range("Q291").Value = '$/bandeja 10 kilos'
$ws.Range("R291").Value = 'Región de O''Higgins'
$ws.Range("S291").Value = 1400
$ws.Range("T291").Value = 10

$ws.Range("D292").Value = 45070
$ws.Range("K292").Value = 'Clementina'
$ws.Range("L292").Value = 'Primera'
$ws.Range("M292").Value = 80
$ws.Range("N292").Value = 12000
$ws.Range("O292").Value = 12000
$ws.Range("P292").Value = 12000
$ws.Range("Q292").Value = '$/bandeja 10 kilos'
$ws.Range("R292").Value = 'Región de O''Higgins'
$ws.Range("S292").Value = 1200
$ws.Range("T292").Value = 10

$ws.Range("D293").Value = 45070
$ws.Range("K293").Value = 'Murcott'
$ws.Range("L293").Value = 'Especial'
$ws.Range("M293").Value = 60
$ws.Range("N293").Value = 14000
$ws.Range("O293").Value = 14000
$ws.Range("P293").Value = 14000
$ws.Range("Q293").Value = '$/bandeja 10 kilos'
$ws.Range("R293").Value = 'Región de O''Higgins'
$ws.Range("S293").Value = 1400
$ws.Range("T293").Value = 10

$ws.Range("D294").Value = 45070
$ws.Range("K294").Value = 'Murcott'
$ws.Range("L294").Value = 'Primera'
$ws.Range("M294").Value = 60
$ws.Range("N294").Value = 12000
$ws.Range("O294").Value = 12000
$ws.Range("P294").Value = 12000
$ws.Range("Q294").Value = '$/bandeja 10 kilos'
$ws.Range("R294").Value = 'Región de O''Higgins'
$ws.Range("S294").Value = 1200
$ws.Range("T294").Value = 10

$ws.Range("D295").Value = 44790
$ws.Range("K295").Value = 'Clementina'
$ws.Range("L295").Value = 'Primera'
$ws.Range("M295").Value = 80
$ws.Range("N295").Value = 8500
$ws.Range("O295").Value = 9000
$ws.Range("P295").Value = 8750
$ws.Range("Q295").Value = '$/caja 18 kilos'
$ws.Range("R295").Value = 'Región de O''Higgins'
$ws.Range("S295").Value = 486
$ws.Range("T295").Value = 18

$ws.Range("D296").Value = 44790
$ws.Range("K296").Value = 'Clementina'
$ws.Range("L296").Value = 'Segunda'
$ws.Range("M296").Value = 120
$ws.Range("N296").Value = 7500
$ws.Range("O296").Value = 8000
$ws.Range("P296").Value = 7750
$ws.Range("Q296").Value = '$/caja 18 kilos'
$ws.Range("R296").Value = 'Región de O''Higgins'
$ws.Range("S296").Value = 431
$ws.Range("T296").Value = 18

$ws.Range("D297").Value = 45040
$ws.Range("K297").Value = 'Murcott'
$ws.Range("L297").Value = 'Especial'
$ws.Range("M297").Value = 80
$ws.Range("N297").Value = 17000
$ws.Range("O297").Value = 17000
$ws.Range("P297").Value = 17000
$ws.Range("Q297").Value = '$/caja 12 kilos'
$ws.Range("R297").Value = 'Región de O''Higgins'
$ws.Range("S297").Value = 1417
$ws.Range("T297").Value = 12

$ws.Range("D298").Value = 45040
$ws.Range("K298").Value = 'Murcott'
$ws.Range("L298").Value = 'Primera'
$ws.Range("M298").Value = 100
$ws.Range("N298").Value = 15000
$ws.Range("O298").Value = 16000
$ws.Range("P298").Value = 15500
$ws.Range("Q298").Value = '$/caja 12 kilos'
$ws.Range("R298").Value = 'Región de O''Higgins'
$ws.Range("S298").Value = 1292
$ws.Range("T298").Value = 12

$ws.Range("D299").Value = 44692
$ws.Range("K299").Value = 'Clemenuless'
$ws.Range("L299").Value = 'Primera'
$ws.Range("M299").Value = 120
$ws.Range("N299").Value = 12000
$ws.Range("O299").Value = 13000
$ws.Range("P299").Value = 12500
$ws.Range("Q299").Value = '$/bandeja 10 kilos'
$ws.Range("R299").Value = 'Provincia de Limarí'
$ws.Range("S299").Value = 1250
$ws.Range("T299").Value = 10

$ws.Range("D300").Value = 45111
$ws.Range("K300").Value = 'Clementina'
$ws.Range("L300").Value = 'Primera'
$ws.Range("M300").Value = 100
$ws.Range("N300").Value = 10000
$ws.Range("O300").Value = 10000
$ws.Range("P300").Value = 10000
$ws.Range("Q300").Value = '$/bandeja 10 kilos'
$ws.Range("R300").Value = 'Región de O''Higgins'
$ws.Range("S300").Value = 1000
$ws.Range("T300").Value = 10

$ws.Range("D301").Value = 45111
$ws.Range("K301").Value = 'Clementina'
$ws.Range("L301").Value = 'Segunda'
$ws.Range("M301").Value = 80
$ws.Range("N301").Value = 8000
$ws.Range("O301").Value = 8000
$ws.Range("P301").Value = 8000
$ws.Range("Q301").Value = '$/bandeja 10 kilos'
$ws.Range("R301").Value = 'Región de O''Higgins'
$ws.Range("S301").Value = 800
$ws.Range("T301").Value = 10

$ws.Range("D302").Value = 44417
$ws.Range("K302").Value = 'Clementina'
$ws.Range("L302").Value = 'Primera'
$ws.Range("M302").Value = 240
$ws.Range("N302").Value = 6000
$ws.Range("O302").Value = 6500
$ws.Range("P302").Value = 6250
$ws.Range("Q302").Value = '$/bandeja 10 kilos'
$ws.Range("R302").Value = 'Región de O''Higgins'
$ws.Range("S302").Value = 625
$ws.Range("T302").Value = 10

$ws.Range("D303").Value = 44417
$ws.Range("K303").Value = 'Clementina'
$ws.Range("L303").Value = 'Segunda'
$ws.Range("M303").Value = 160
$ws.Range("N303").Value = 5000
$ws.Range("O303").Value = 5500
$ws.Range("P303").Value = 5250
$ws.Range("Q303").Value = '$/bandeja 10 kilos'
$ws.Range("R303").Value = 'Región de O''Higgins'
$ws.Range("S303").Value = 525
$ws.Range("T303").Value = 10

$ws.Range("D304").Value = 44323
$ws.Range("K304").Value = 'Murcott'
$ws.Range("L304").Value = 'Primera'
$ws.Range("M304").Value = 120
$ws.Range("N304").Value = 13000
$ws.Range("O304").Value = 14000
$ws.Range("P304").Value = 13500
$ws.Range("Q304").Value = '$/bandeja 10 kilos'
$ws.Range("R304").Value = 'Provincia de Limarí'
$ws.Range("S304").Value = 1350
$ws.Range("T304").Value = 10

$ws.Range("D305").Value = 44468
$ws.Range("K305").Value = 'Murcott'
$ws.Range("L305").Value = 'Primera'
$ws.Range("M305").Value = 300
$ws.Range("N305").Value = 6000
$ws.Range("O305").Value = 6500
$ws.Range("P305").Value = 6250
$ws.Range("Q305").Value = '$/bandeja 10 kilos'
$ws.Range("R305").Value = 'Provincia de Limarí'
$ws.Range("S305").Value = 625
$ws.Range("T305").Value = 10

$ws.Range("D306").Value = 44468
$ws.Range("K306").Value = 'Murcott'
$ws.Range("L306").Value = 'Segunda'
$ws.Range("M306").Value = 120
$ws.Range("N306").Value = 5500
$ws.Range("O306").Value = 5500
$ws.Range("P306").Value = 5500
$ws.Range("Q306").Value = '$/bandeja 10 kilos'
$ws.Range("R306").Value = 'Provincia de Limarí'
$ws.Range("S306").Value = 550
$ws.Range("T306").Value = 10

$ws.Range("D307").Value = 45106
$ws.Range("K307").Value = 'Clementina'
$ws.Range("L307").Value = 'Especial'
$ws.Range("M307").Value = 60
$ws.Range("N307").Value = 11000
$ws.Range("O307").Value = 11000
$ws.Range("P307").Value = 11000
$ws.Range("Q307").Value = '$/bandeja 10 kilos'
$ws.Range("R307").Value = 'Región de O''Higgins'
$ws.Range("S307").Value = 1100
$ws.Range("T307").Value = 10

$ws.Range("D308").Value = 45106
$ws.Range("K308").Value = 'Clementina'
$ws.Range("L308").Value = 'Primera'
$ws.Range("M308").Value = 50
$ws.Range("N308").Value = 9000
$ws.Range("O308").Value = 9000
$ws.Range("P308").Value = 9000
$ws.Range("Q308").Value = '$/bandeja 10 kilos'
$ws.Range("R308").Value = 'Región de O''Higgins'
$ws.Range("S308").Value = 900
$ws.Range("T308").Value = 10

$ws.Range("D309").Value = 45106
$ws.Range("K309").Value = 'Clementina'
$ws.Range("L309").Value = 'Segunda'
$ws.Range("M309").Value = 50
$ws.Range("N309").Value = 8000
$ws.Range("O309").Value = 8000
$ws.Range("P309").Value = 8000
$ws.Range("Q309").Value = '$/bandeja 10 kilos'
$ws.Range("R309").Value = 'Región de O''Higgins'
$ws.Range("S309").Value = 800
$ws.Range("T309").Value = 10

$ws.Range("D310").Value = 44358
$ws.Range("K310").Value = 'Clemenuless'
$ws.Range("L310").Value = 'Primera'
$ws.Range("M310").Value = 60
$ws.Range("N310").Value = 11000
$ws.Range("O310").Value = 11000
$ws.Range("P310").Value = 11000
$ws.Range("Q310").Value = '$/bandeja 10 kilos'
$ws.Range("R310").Value = 'Provincia de Limarí'
$ws.Range("S310").Value = 1100
$ws.Range("T310").Value = 10

$ws.Range("D311").Value = 44358
$ws.Range("K311").Value = 'Clemenuless'
$ws.Range("L311").Value = 'Segunda'
$ws.Range("M311").Value = 120
$ws.Range("N311").Value = 9000
$ws.Range("O311").Value = 10000
$ws.Range("P311").Value = 9500
$ws.Range("Q311").Value = '$/bandeja 10 kilos'
$ws.Range("R311").Value = 'Provincia de Limarí'
$ws.Range("S311").Value = 950
$ws.Range("T311").Value = 10

$ws.Range("D312").Value = 44391
$ws.Range("K312").Value = 'Clemenuless'
$ws.Range("L312").Value = 'Primera'
$ws.Range("M312").Value = 120
$ws.Range("N312").Value = 7000
$ws.Range("O312").Value = 7500
$ws.Range("P312").Value = 7250
$ws.Range("Q312").Value = '$/bandeja 10 kilos'
$ws.Range("R312").Value = 'Provincia de Limarí'
$ws.Range("S312").Value = 725
$ws.Range("T312").Value = 10

$ws.Range("D313").Value = 44391
$ws.Range("K313").Value = 'Clemenuless'
$ws.Range("L313").Value = 'Segunda'
$ws.Range("M313").Value = 120
$ws.Range("N313").Value = 6000
$ws.Range("O313").Value = 6500
$ws.Range("P313").Value = 6250
$ws.Range("Q313").Value = '$/bandeja 10 kilos'
$ws.Range("R313").Value = 'Provincia de Limarí'
$ws.Range("S313").Value = 625
$ws.Range("T313").Value = 10

$ws.Range("D314").Value = 44420
$ws.Range("K314").Value = 'Clementina'
$ws.Range("L314").Value = 'Primera'
$ws.Range("M314").Value = 200
$ws.Range("N314").Value = 6000
$ws.Range("O314").Value = 6500
$ws.Range("P314").Value = 6250
$ws.Range("Q314").Value = '$/bandeja 10 kilos'
$ws.Range("R314").Value = 'Región de O''Higgins'
$ws.Range("S314").Value = 625
$ws.Range("T314").Value = 10

$ws.Range("D315").Value = 44420
$ws.Range("K315").Value = 'Clementina'
$ws.Range("L315").Value = 'Segunda'
$ws.Range("M315").Value = 200
$ws.Range("N315").Value = 5000
$ws.Range("O315").Value = 5500
$ws.Range("P315").Value = 5250
$ws.Range("Q315").Value = '$/bandeja 10 kilos'
$ws.Range("R315").Value = 'Región de O''Higgins'
$ws.Range("S315").Value = 525
$ws.Range("T315").Value = 10

$ws.Range("D316").Value = 45075
$ws.Range("K316").Value = 'Clementina'
$ws.Range("L316").Value = 'Especial'
$ws.Range("M316").Value = 40
$ws.Range("N316").Value = 14000
$ws.Range("O316").Value = 14000
$ws.Range("P316").Value = 14000
$ws.Range("Q316").Value = '$/bandeja 10 kilos'
$ws.Range("R316").Value = 'Región de O''Higgins'
$ws.Range("S316").Value = 1400
$ws.Range("T316").Value = 10

$ws.Range("D317").Value = 45075
$ws.Range("K317").Value = 'Clementina'
$ws.Range("L317").Value = 'Primera'
$ws.Range("M317").Value = 50
$ws.Range("N317").Value = 12000
$ws.Range("O317").Value = 12000
$ws.Range("P317").Value = 12000
$ws.Range("Q317").Value = '$/bandeja 10 kilos'
$ws.Range("R317").Value = 'Región de O''Higgins'
$ws.Range("S317").Value = 1200
$ws.Range("T317").Value = 10

$ws.Range("D318").Value = 45075
$ws.Range("K318").Value = 'Clementina'
$ws.Range("L318").Value = 'Segunda'
$ws.Range("M318").Value = 30
$ws.Range("N318").Value = 10000
$ws.Range("O318").Value = 10000
$ws.Range("P318").Value = 10000
$ws.Range("Q318").Value = '$/bandeja 10 kilos'
$ws.Range("R318").Value = 'Región de O''Higgins'
$ws.Range("S318").Value = 1000
$ws.Range("T318").Value = 10

$ws.Range("D319").Value = 44364
$ws.Range("K319").Value = 'Clementina'
$ws.Range("L319").Value = 'Primera'
$ws.Range("M319").Value = 120
$ws.Range("N319").Value = 8000
$ws.Range("O319").Value = 8500
$ws.Range("P319").Value = 8250
$ws.Range("Q319").Value = '$/bandeja 10 kilos'
$ws.Range("R319").Value = 'Provincia de Limarí'
$ws.Range("S319").Value = 825
$ws.Range("T319").Value = 10

$ws.Range("D320").Value = 44364
$ws.Range("K320").Value = 'Clementina'
$ws.Range("L320").Value = 'Segunda'
$ws.Range("M320").Value = 120
$ws.Range("N320").Value = 7000
$ws.Range("O320").Value = 7500
$ws.Range("P320").Value = 7250
$ws.Range("Q320").Value = '$/bandeja 10 kilos'
$ws.Range("R320").Value = 'Provincia de Limarí'
$ws.Range("S320").Value = 725
$ws.Range("T320").Value = 10

$ws.Range("D321").Value = 45112
$ws.Range("K321").Value = 'Clementina'
$ws.Range("L321").Value = 'Primera'
$ws.Range("M321").Value = 80
$ws.Range("N321").Value = 10000
$ws.Range("O321").Value = 10000
$ws.Range("P321").Value = 10000
$ws.Range("Q321").Value = '$/bandeja 10 kilos'
$ws.Range("R321").Value = 'Región de O''Higgins'
$ws.Range("S321").Value = 1000
$ws.Range("T321").Value = 10

$ws.Range("D322").Value = 45112
$ws.Range("K322").Value = 'Clementina'
$ws.Range("L322").Value = 'Segunda'
$ws.Range("M322").Value = 80
$ws.Range("N322").Value = 8000
$ws.Range("O322").Value = 8000
$ws.Range("P322").Value = 8000
$ws.Range("Q322").Value = '$/bandeja 10 kilos'
$ws.Range("R322").Value = 'Región de O''Higgins'
$ws.Range("S322").Value = 800
$ws.Range("T322").Value = 10

$ws.Range("D323").Value = 45112
$ws.Range("K323").Value = 'Clementina'
$ws.Range("L323").Value = 'Primera'
$ws.Range("M323").Value = 100
$ws.Range("N323").Value = 10000
$ws.Range("O323").Value = 10000
$ws.Range("P323").Value = 10000
$ws.Range("Q323").Value = '$/bandeja 10 kilos'
$ws.Range("R323").Value = 'Región de O''Higgins'
$ws.Range("S323").Value = 1000
$ws.Range("T323").Value = 10

$ws.Range("D324").Value = 45112
$ws.Range("K324").Value = 'Clementina'
$ws.Range("L324").Value = 'Segunda'
$ws.Range("M324").Value = 80
$ws.Range("N324").Value = 8000
$ws.Range("O324").Value = 8000
$ws.Range("P324").Value = 8000
$ws.Range("Q324").Value = '$/bandeja 10 kilos'
$ws.Range("R324").Value = 'Región de O''Higgins'
$ws.Range("S324").Value = 800
$ws.Range("T324").Value = 10

$ws.Range("D325").Value = 44802
$ws.Range("K325").Value = 'Murcott'
$ws.Range("L325").Value = 'Primera'
$ws.Range("M325").Value = 120
$ws.Range("N325").Value = 8000
$ws.Range("O325").Value = 8500
$ws.Range("P325").Value = 8250
$ws.Range("Q325").Value = '$/caja 18 kilos'
$ws.Range("R325").Value = 'Región de O''Higgins'
$ws.Range("S325").Value = 458
$ws.Range("T325").Value = 18

$ws.Range("D326").Value = 44802
$ws.Range("K326").Value = 'Murcott'
$ws.Range("L326").Value = 'Segunda'
$ws.Range("M326").Value = 60
$ws.Range("N326").Value = 7000
$ws.Range("O326").Value = 7000
$ws.Range("P326").Value = 7000
$ws.Range("Q326").Value = '$/caja 18 kilos'
$ws.Range("R326").Value = 'Región de O''Higgins'
$ws.Range("S326").Value = 389
$ws.Range("T326").Value = 18

$ws.Range("D327").Value = 44348
$ws.Range("K327").Value = 'Clemenuless'
$ws.Range("L327").Value = 'Primera'
$ws.Range("M327").Value = 120
$ws.Range("N327").Value = 10000
$ws.Range("O327").Value = 11000
$ws.Range("P327").Value = 10500
$ws.Range("Q327").Value = '$/bandeja 10 kilos'
$ws.Range("R327").Value = 'Provincia de Limarí'
$ws.Range("S327").Value = 1050
$ws.Range("T327").Value = 10

$ws.Range("D328").Value = 45089
$ws.Range("K328").Value = 'Clementina'
$ws.Range("L328").Value = 'Especial'
$ws.Range("M328").Value = 80
$ws.Range("N328").Value = 12000
$ws.Range("O328").Value = 12000
$ws.Range("P328").Value = 12000
$ws.Range("Q328").Value = '$/bandeja 10 kilos'
$ws.Range("R328").Value = 'Región de O''Higgins'
$ws.Range("S328").Value = 1200
$ws.Range("T328").Value = 10

$ws.Range("D329").Value = 45089
$ws.Range("K329").Value = 'Clementina'
$ws.Range("L329").Value = 'Primera'
$ws.Range("M329").Value = 60
$ws.Range("N329").Value = 10000
$ws.Range("O329").Value = 10000
$ws.Range("P329").Value = 10000
$ws.Range("Q329").Value = '$/bandeja 10 kilos'
$ws.Range("R329").Value = 'Región de O''Higgins'
$ws.Range("S329").Value = 1000
$ws.Range("T329").Value = 10

$ws.Range("D330").Value = 45089
$ws.Range("K330").Value = 'Clementina'
$ws.Range("L330").Value = 'Segunda'
$ws.Range("M330").Value = 50
$ws.Range("N330").Value = 8000
$ws.Range("O330").Value = 8000
$ws.Range("P330").Value = 8000
$ws.Range("Q330").Value = '$/bandeja 10 kilos'
$ws.Range("R330").Value = 'Región de O''Higgins'
$ws.Range("S330").Value = 800
$ws.Range("T330").Value = 10

$ws.Range("D331").Value = 45099
$ws.Range("K331").Value = 'Clementina'
$ws.Range("L331").Value = 'Especial'
$ws.Range("M331").Value = 80
$ws.Range("N331").Value = 12000
$ws.Range("O331").Value = 12000
$ws.Range("P331").Value = 12000
$ws.Range("Q331").Value = '$/bandeja 10 kilos'
$ws.Range("R331").Value = 'Región de O''Higgins'
$ws.Range("S331").Value = 1200
$ws.Range("T331").Value = 10

$ws.Range("D332").Value = 45099
$ws.Range("K332").Value = 'Clementina'
$ws.Range("L332").Value = 'Primera'
$ws.Range("M332").Value = 60
$ws.Range("N332").Value = 10000
$ws.Range("O332").Value = 10000
$ws.Range("P332").Value = 10000
$ws.Range("Q332").Value = '$/bandeja 10 kilos'
$ws.Range("R332").Value = 'Región de O''Higgins'
$ws.Range("S332").Value = 1000
$ws.Range("T332").Value = 10

$ws.Range("D333").Value = 45099
$ws.Range("K333").Value = 'Clementina'
$ws.Range("L333").Value = 'Segunda'
$ws.Range("M333").Value = 50
$ws.Range("N333").Value = 8000
$ws.Range("O333").Value = 8000
$ws.Range("P333").Value = 8000
$ws.Range("Q333").Value = '$/bandeja 10 kilos'
$ws.Range("R333").Value = 'Región de O''Higgins'
$ws.Range("S333").Value = 800
$ws.Range("T333").Value = 10

$ws.Range("D334").Value = 45121
$ws.Range("K334").Value = 'Clementina'
$ws.Range("L334").Value = 'Primera'
$ws.Range("M334").Value = 80
$ws.Range("N334").Value = 10000
$ws.Range("O334").Value = 10000
$ws.Range("P334").Value = 10000
$ws.Range("Q334").Value = '$/bandeja 10 kilos'
$ws.Range("R334").Value = 'Región de O''Higgins'
$ws.Range("S334").Value = 1000
$ws.Range("T334").Value = 10

$ws.Range("D335").Value = 45121
$ws.Range("K335").Value = 'Clementina'
$ws.Range("L335").Value = 'Segunda'
$ws.Range("M335").Value = 60
$ws.Range("N335").Value = 8000
$ws.Range("O335").Value = 8000
$ws.Range("P335").Value = 8000
$ws.Range("Q335").Value = '$/bandeja 10 kilos'
$ws.Range("R335").Value = 'Región de O''Higgins'
$ws.Range("S335").Value = 800
$ws.Range("T335").Value = 10

$ws.Range("D336").Value = 44341
$ws.Range("K336").Value = 'Clemenuless'
$ws.Range("L336").Value = 'Primera'
$ws.Range("M336").Value = 120
$ws.Range("N336").Value = 11000
$ws.Range("O336").Value = 12000
$ws.Range("P336").Value = 11500
$ws.Range("Q336").Value = '$/bandeja 10 kilos'
$ws.Range("R336").Value = 'Provincia de Limarí'
$ws.Range("S336").Value = 1150
$ws.Range("T336").Value = 10

$ws.Range("D337").Value = 44341
$ws.Range("K337").Value = 'Clemenuless'
$ws.Range("L337").Value = 'Segunda'
$ws.Range("M337").Value = 80
$ws.Range("N337").Value = 10000
$ws.Range("O337").Value = 10000
$ws.Range("P337").Value = 10000
$ws.Range("Q337").Value = '$/bandeja 10 kilos'
$ws.Range("R337").Value = 'Provincia de Limarí'
$ws.Range("S337").Value = 1000
$ws.Range("T337").Value = 10

$ws.Range("D338").Value = 45072
$ws.Range("K338").Value = 'Murcott'
$ws.Range("L338").Value = 'Especial'
$ws.Range("M338").Value = 50
$ws.Range("N338").Value = 14000
$ws.Range("O338").Value = 14000
$ws.Range("P338").Value = 14000
$ws.Range("Q338").Value = '$/bandeja 10 kilos'
$ws.Range("R338").Value = 'Región de O''Higgins'
$ws.Range("S338").Value = 1400
$ws.Range("T338").Value = 10

$ws.Range("A339").Value = 7
$ws.Range("B339").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C339").Value = 'Ñuble'
$ws.Range("E339").Value = 16
$ws.Range("F339").Value = 'Fruta'
$ws.Range("G339").Value = 100102
$ws.Range("H339").Value = 'Cítricos'
$ws.Range("I339").Value = 100102004
$ws.Range("J339").Value = 'Mandarina'
$ws.Range("D339").Value = 45072
$ws.Range("K339").Value = 'Murcott'
$ws.Range("L339").Value = 'Primera'
$ws.Range("M339").Value = 80
$ws.Range("N339").Value = 12000
$ws.Range("O339").Value = 12000
$ws.Range("P339").Value = 12000
$ws.Range("Q339").Value = '$/bandeja 10 kilos'
$ws.Range("R339").Value = 'Región de O''Higgins'
$ws.Range("S339").Value = 1200
$ws.Range("T339").Value = 10
$ws.Range("D339").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A340").Value = 7
$ws.Range("B340").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C340").Value = 'Ñuble'
$ws.Range("E340").Value = 16
$ws.Range("F340").Value = 'Fruta'
$ws.Range("G340").Value = 100102
$ws.Range("H340").Value = 'Cítricos'
$ws.Range("I340").Value = 100102004
$ws.Range("J340").Value = 'Mandarina'
$ws.Range("D340").Value = 45072
$ws.Range("K340").Value = 'Murcott'
$ws.Range("L340").Value = 'Segunda'
$ws.Range("M340").Value = 80
$ws.Range("N340").Value = 10000
$ws.Range("O340").Value = 10000
$ws.Range("P340").Value = 10000
$ws.Range("Q340").Value = '$/bandeja 10 kilos'
$ws.Range("R340").Value = 'Región de O''Higgins'
$ws.Range("S340").Value = 1000
$ws.Range("T340").Value = 10
$ws.Range("D340").NumberFormat = "YYYY-MM-DD HH:MM:SS"
